$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "정상적으로 로딩되었습니다.`n좋은 하루 보내세요!"
$ws.Range("I3").Value = "상황실"
$ws.Range("I4").Value = "최신 버전 업데이트가 있습니다."
$ws.Range("I5").Value = "버전 업데이트가 가능합니다!"
$ws.Range("I6").Value = "현재 버전`n최신 버전`n"
$ws.Range("I8").Value = "상황종료- 코드 4"
$ws.Range("P8").Value = "Jste ~g~pod kodem 4~s~.`nNeni potreba zadna dalsi jednotka."
$ws.Range("I10").Value = "가까이 오세요!"
$ws.Range("P10").Value = "Prilis daleko.`nProsim, priblizte se."
$ws.Range("I12").Value = "남성"
$ws.Range("P12").Value = "Muz"
$ws.Range("I13").Value = "여성"
$ws.Range("P13").Value = "Zena"
$ws.Range("P15").Value = "Automaticka aktualizace"
$ws.Range("P16").Value = "~r~Automaticka aktualizace selhala.~s~`nAktualizujte prosim ~y~rucne~ś~."
$ws.Range("P17").Value = "~y~Aktualizuji~s~ ~b~{0}~s~ nyni…"
$ws.Range("P18").Value = "~b~{0}~s~ byl ~r~automaticky aktualizovan~s~.`nProsim ~y~restartujte~s~ ~b~LSPDFR~s~."
